# Auto commit at 2025-08-20  7:15:02.04
#
# Refresh the day's figures on the "Metrics" sheet and mirror the same
# numbers into the "today" sheet's helper columns (B/E), then leave the
# selection / active sheet exactly where the author left them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Metrics sheet: B2:B13 updated with the new cumulative figures.
# ---------------------------------------------------------------------
$wsMetrics = $wb.Worksheets.Item("Metrics")

$metricsValues = @{
    "B2"  = 320961.33
    "B3"  = 275790.73
    "B4"  = 101191.54000000001
    "B5"  = 12649
    "B6"  = 3716589.9
    "B7"  = 3154505.39
    "B8"  = 1064834.0999999999
    "B9"  = 143337
    "B10" = 32181913.699999999
    "B11" = 19184375.460000001
    "B12" = 11346542.99
    "B13" = 1240964
}

foreach ($addr in $metricsValues.Keys) {
    $wsMetrics.Range($addr).Value = $metricsValues[$addr]
}

$null = $wsMetrics.Activate()
$null = $wsMetrics.Range("B2:B13").Select()

# ---------------------------------------------------------------------
# today sheet: same underlying numbers mirrored into B11:B22 / E11:E22,
# the now-empty B6 cell gets the shared "#,##0.00" style (index 4)
# instead of its previous one-off style, and the selection moves to J16.
# ---------------------------------------------------------------------
$wsToday = $wb.Worksheets.Item("today")

$todayValues = @{
    11 = 320961.33
    12 = 275790.73
    13 = 101191.54000000001
    14 = 12649
    15 = 3716589.9
    16 = 3154505.39
    17 = 1064834.0999999999
    18 = 143337
    19 = 32181913.699999999
    20 = 19184375.460000001
    21 = 11346542.99
    22 = 1240964
}

foreach ($row in $todayValues.Keys) {
    $v = $todayValues[$row]
    $wsToday.Cells.Item($row, 2).Value = $v   # column B
    $wsToday.Cells.Item($row, 5).Value = $v   # column E
}

# B6 keeps its same (blank) content but switches to the shared style
# used by B3:B5 (#,##0.00), which drops the now-unused one-off xf.
$wsToday.Cells.Item(6, 2).NumberFormat = $wsToday.Cells.Item(3, 2).NumberFormat

$null = $wsToday.Activate()
$null = $wsToday.Range("J16").Select()
